# Apply scheduled-runner market-data updates to the Leve profit sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 27780750
$ws.Range("I113").Value = 71430690
$ws.Range("J113").Value = 3516.182
$ws.Range("K113").Value = 71430690
$ws.Range("L113").Value = 3516.182
$ws.Range("M113").Value = -71427436
$ws.Range("N113").Value = -10024.182

$ws.Range("H137").Value = 1339.6957
$ws.Range("I137").Value = 1174.375
$ws.Range("J137").Value = 1717.5714
$ws.Range("K137").Value = 3523.125
$ws.Range("L137").Value = 5152.7142
$ws.Range("M137").Value = -973.125
$ws.Range("N137").Value = -10252.7142

$ws.Range("H138").Value = 2472.8975
$ws.Range("I138").Value = 1627.8077
$ws.Range("J138").Value = 4163.077
$ws.Range("K138").Value = 4883.4231
$ws.Range("L138").Value = 12489.231
$ws.Range("M138").Value = 256.5769
$ws.Range("N138").Value = -22769.231

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H76").Value = 121361
$ws.Range("J76").Value = 121361
$ws.Range("L76").Value = 121361
$ws.Range("N76").Value = -122037

$ws.Range("H79").Value = 121361
$ws.Range("J79").Value = 121361
$ws.Range("L79").Value = 121361
$ws.Range("N79").Value = -123701

$ws.Range("H88").Value = 4455.3335
$ws.Range("I88").Value = 4500
$ws.Range("J88").Value = 4449.75
$ws.Range("K88").Value = 4500
$ws.Range("L88").Value = 4449.75
$ws.Range("M88").Value = -4094
$ws.Range("N88").Value = -5261.75

$ws.Range("H91").Value = 4455.3335
$ws.Range("I91").Value = 4500
$ws.Range("J91").Value = 4449.75
$ws.Range("K91").Value = 4500
$ws.Range("L91").Value = 4449.75
$ws.Range("M91").Value = -3096
$ws.Range("N91").Value = -7257.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 464.08334
$ws.Range("J64").Value = 629.5
$ws.Range("L64").Value = 629.5
$ws.Range("N64").Value = -1079.5

$ws.Range("H67").Value = 464.08334
$ws.Range("J67").Value = 629.5
$ws.Range("L67").Value = 629.5
$ws.Range("N67").Value = -2189.5

$ws.Range("H86").Value = 2459.862
$ws.Range("I86").Value = 2580.375
$ws.Range("J86").Value = 1881.4
$ws.Range("K86").Value = 2580.375
$ws.Range("L86").Value = 1881.4
$ws.Range("M86").Value = -1457.375
$ws.Range("N86").Value = -4127.4

$ws.Range("H89").Value = 2459.862
$ws.Range("I89").Value = 2580.375
$ws.Range("J89").Value = 1881.4
$ws.Range("K89").Value = 12901.875
$ws.Range("L89").Value = 9407
$ws.Range("M89").Value = -7285.875
$ws.Range("N89").Value = -20639

$ws.Range("H99").Value = 3316.7778
$ws.Range("I99").Value = 2047.5
$ws.Range("J99").Value = 4332.2
$ws.Range("K99").Value = 2047.5
$ws.Range("L99").Value = 4332.2
$ws.Range("M99").Value = -549.5
$ws.Range("N99").Value = -7328.2

$ws.Range("H107").Value = 1601.4667
$ws.Range("I107").Value = 1460.2858
$ws.Range("J107").Value = 1725
$ws.Range("K107").Value = 1460.2858
$ws.Range("L107").Value = 1725
$ws.Range("M107").Value = 459.7141999999999
$ws.Range("N107").Value = -5565

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2043.963
$ws.Range("I58").Value = 2309.158
$ws.Range("J58").Value = 1414.125
$ws.Range("K58").Value = 2309.158
$ws.Range("L58").Value = 1414.125
$ws.Range("M58").Value = -2106.158
$ws.Range("N58").Value = -1820.125

$ws.Range("H86").Value = 5872.143
$ws.Range("I86").Value = 8623.182000000001
$ws.Range("J86").Value = 2846
$ws.Range("K86").Value = 8623.182000000001
$ws.Range("L86").Value = 2846
$ws.Range("M86").Value = -7500.182000000001
$ws.Range("N86").Value = -5092

$ws.Range("H89").Value = 5872.143
$ws.Range("I89").Value = 8623.182000000001
$ws.Range("J89").Value = 2846
$ws.Range("K89").Value = 43115.91
$ws.Range("L89").Value = 14230
$ws.Range("M89").Value = -37499.91
$ws.Range("N89").Value = -25462

$ws.Range("H107").Value = 1210.6364
$ws.Range("I107").Value = 1281.3
$ws.Range("J107").Value = 1101.9231
$ws.Range("K107").Value = 1281.3
$ws.Range("L107").Value = 1101.9231
$ws.Range("M107").Value = 638.7
$ws.Range("N107").Value = -4941.9231

$ws.Range("H133").Value = 31151.5
$ws.Range("J133").Value = 31151.5
$ws.Range("L133").Value = 31151.5
$ws.Range("N133").Value = -36211.5

$ws.Range("H136").Value = 2043.963
$ws.Range("I136").Value = 2309.158
$ws.Range("J136").Value = 1414.125
$ws.Range("K136").Value = 6927.474
$ws.Range("L136").Value = 4242.375
$ws.Range("M136").Value = -4377.474
$ws.Range("N136").Value = -9342.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 448.56097
$ws.Range("I113").Value = 394.8095
$ws.Range("J113").Value = 505
$ws.Range("K113").Value = 1184.4285
$ws.Range("L113").Value = 1515
$ws.Range("M113").Value = 985.5715
$ws.Range("N113").Value = -5855

$ws.Range("H122").Value = 3103.1965
$ws.Range("I122").Value = 332.55554
$ws.Range("J122").Value = 3633.7446
$ws.Range("K122").Value = 2992.99986
$ws.Range("L122").Value = 32703.7014
$ws.Range("M122").Value = -542.9998599999999
$ws.Range("N122").Value = -37603.7014

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 8480
$ws.Range("J5").Value = 9960
$ws.Range("L5").Value = 9960
$ws.Range("N5").Value = -10184

$ws.Range("H45").Value = 13420.3
$ws.Range("J45").Value = 13420.3
$ws.Range("L45").Value = 13420.3
$ws.Range("N45").Value = -14538.3

$ws.Range("H107").Value = 669.62964
$ws.Range("I107").Value = 459.42105
$ws.Range("J107").Value = 1168.875
$ws.Range("K107").Value = 459.42105
$ws.Range("L107").Value = 1168.875
$ws.Range("M107").Value = 1460.57895
$ws.Range("N107").Value = -5008.875

$ws.Range("H113").Value = 2000
$ws.Range("I113").Value = 2000
$ws.Range("K113").Value = 2000
$ws.Range("M113").Value = 170

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2133.3333
$ws.Range("J100").Value = 3000
$ws.Range("L100").Value = 3000
$ws.Range("N100").Value = -4082

$ws.Range("H132").Value = 1745.5333
$ws.Range("I132").Value = 1384.6818
$ws.Range("J132").Value = 2737.875
$ws.Range("K132").Value = 4154.0454
$ws.Range("L132").Value = 8213.625
$ws.Range("M132").Value = -1624.0454
$ws.Range("N132").Value = -13273.625
